$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells keep their original Text format so that
# numeric-looking strings (e.g. "0.991") are not auto-converted to numbers.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.129.38"
$ws.Range("E2").Value = "  +1.44%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.641.80"
$ws.Range("E3").Value = "  +0.53%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.991"
$ws.Range("E4").Value = "  -1.01%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.61"
$ws.Range("E5").Value = "  +0.33%  "

$ws.Range("E6").Value = "  +0.01%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.993"
$ws.Range("E7").Value = "  -0.82%  "

$ws.Range("E8").Value = "  -0.23%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0633"
$ws.Range("E9").Value = "  -0.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.84"
$ws.Range("E10").Value = "  +0.88%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0788"
$ws.Range("E11").Value = "  +0.29%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.713.16"
$ws.Range("E12").Value = "  +4.91%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").Value = "4.26"
$ws.Range("E13").Value = "  +0.36%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.871.41"
$ws.Range("E14").Value = "  +0.75%  "

$ws.Range("E15").Value = "  -0.39%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₃0765"
$ws.Range("E16").Value = "  +0.19%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.47"
$ws.Range("E17").Value = "  +1.22%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.131.26"
$ws.Range("E18").Value = "  +1.47%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.993"
$ws.Range("E19").Value = "  -0.83%  "

$ws.Range("E20").Value = "  +0.54%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "193.55"
$ws.Range("E21").Value = "  -0.26%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.04"
$ws.Range("E22").Value = "  +1.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.42"
$ws.Range("E23").Value = "  +2.54%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.992"
$ws.Range("E24").Value = "  -0.97%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "141.96"
$ws.Range("E26").Value = "  -0.63%  "

$ws.Range("E27").Value = "  +0.91%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.91"
$ws.Range("E28").Value = "  +0.57%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.67"
$ws.Range("E29").Value = "  +0.82%  "

$ws.Range("E30").Value = "  +0.06%  "

$ws.Range("E31").Value = "  +0.51%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.34"
$ws.Range("E32").Value = "  +0.57%  "

$ws.Range("E33").Value = "  +0.16%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.61"
$ws.Range("E34").Value = "  +1.74%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.38"
$ws.Range("E35").Value = "  -0.40%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.911"
$ws.Range("E36").Value = "  +1.13%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.153.75"
$ws.Range("E37").Value = "  +2.23%  "

$ws.Range("E38").Value = "  +0.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.49"
$ws.Range("E39").Value = "  -1.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.993"
$ws.Range("E41").Value = "  -0.92%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.66"
$ws.Range("E42").Value = "  +1.61%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.70"
$ws.Range("E43").Value = "  +0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.804"
$ws.Range("E44").Value = "  -0.14%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.781.52"
$ws.Range("E45").Value = "  +0.80%  "

$ws.Range("E46").Value = "  -2.60%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.83"
$ws.Range("E47").Value = "  +1.42%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("B48").Value = "Cronos"
$ws.Range("C48").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D48").Value = "0.0519"
$ws.Range("E48").Value = "  +3.10%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("B49").Value = "RenderToken"
$ws.Range("C49").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D49").Value = "1.47"
$ws.Range("E49").Value = "  +6.34%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("B50").Value = "Mantle"
$ws.Range("C50").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D50").Value = "0.416"
$ws.Range("E50").Value = "  -0.10%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "7.66"
$ws.Range("E51").Value = "  +1.24%  "

